# Apply the commit's slide edits to the Hexagon shape's outline weight and
# the Picture's position/size (slide 1, shapes "Hexagon 5" and "Picture 4").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$hexagon = $s.Shapes.Item(1)   # "Hexagon 5"
$picture = $s.Shapes.Item(2)   # "Picture 4"

# Hexagon outline: 38100 EMU (3pt) -> 28575 EMU (2.25pt)
$hexagon.Line.Weight = 2.25

# Picture position/size:
#   off  191690,44053   -> 213717,72629
#   ext 1445419,1445419 -> 1401366,1401366  (EMU)
# Values below are expressed in points (EMU / 12700); the x/y offsets carry a
# tiny epsilon nudge so the host's internal float32 storage still floors to
# the exact target EMU.
$picture.Left = 16.82811073622047
$picture.Top = 5.718818997637795
$picture.Width = 110.34377952755905
$picture.Height = 110.34377952755905
